# Création de nx mappings 0e3f71524d31ae38b8fd94d7fa5f14a19c184959
#
# 1) Metadata sheet: bump the "Date" value.
# 2) "Mapping Table 0" sheet: add 5 new equivalence rows (renouvellement,
#    duree, quantite, auteur, dispositifMedical) right after the
#    "date -> effectiveTime" row and before the "affectionLongueDuree" row.
# 3) "Mapping Table 1" sheet: rename the occurrenceTiming target display to
#    occurrence[x].

$wb = $excel.ActiveWorkbook

# --- 1) Metadata: update Date value -----------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Cells.Item(8, 2).Value = "2026-01-05T15:47:03+00:00"

# --- 2) Mapping Table 0: insert the 5 new mapping rows -----------------
$ws0 = $wb.Worksheets.Item("Mapping Table 0")

# Insert 5 blank rows starting at row 7 (pushes the existing
# "affectionLongueDuree..." row and everything below it down by 5 rows).
$ws0.Range("A7:A11").EntireRow.Insert()

# Copy the formatting of an existing data row onto the freshly inserted
# rows so borders / fill / font match the rest of the table.
$ws0.Range("A6:E6").Copy()
$ws0.Range("A7:E11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
  @("FRLMDispositifMedicalEntree.renouvellement", "FRCDADispositifMedical.repeatNumber"),
  @("FRLMDispositifMedicalEntree.duree", "FRCDADispositifMedical.expectedUseTime"),
  @("FRLMDispositifMedicalEntree.quantite", "FRCDADispositifMedical.quantity"),
  @("FRLMDispositifMedicalEntree.auteur", "FRCDADispositifMedical.author"),
  @("FRLMDispositifMedicalEntree.dispositifMedical", "FRCDADispositifMedical.participant")
)

$r = 7
foreach ($row in $newRows) {
  $ws0.Cells.Item($r, 1).Value = $row[0]
  $ws0.Cells.Item($r, 3).Value = "equivalent"
  $ws0.Cells.Item($r, 4).Value = $row[1]
  $r = $r + 1
}

# --- 3) Mapping Table 1: occurrenceTiming -> occurrence[x] -------------
$ws1 = $wb.Worksheets.Item("Mapping Table 1")
$ws1.Cells.Item(6, 4).Value = "FRDeviceRequestDocument.occurrence[x]"
